$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 11-14 (row 11 previously only had A filled in)
$ws.Range("B11").Value = "OptionText: TitleLanguage"
$ws.Range("C11").Value = "Language"
$ws.Range("D11").Value = "Lenguaje"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "OptionButton: BackButtonText"
$ws.Range("C12").Value = "Back"
$ws.Range("D12").Value = "Atrás"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "OptionButton: HomeButtonText"
$ws.Range("C13").Value = "Home"
$ws.Range("D13").Value = "Inicio"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "OptionButton: TitleResolutionMute"
$ws.Range("C14").Value = "FullScreen"
$ws.Range("D14").Value = "Pantalla Completa"

# D1: "Spanish" -> "Español" (set last so it becomes the last unique shared string)
$ws.Range("D1").Value = "Español"

# The previous styled-but-empty cells (C12, D12, B13) are replaced by plain (unstyled) values above,
# so clear any leftover style by resetting to the default "Normal" style.
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("B13").Style = "Normal"

# Update the active selection to D11
$ws.Range("D11").Select()
